$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1070.5714
$ws.Range("I33").Value = 173.75
$ws.Range("J33").Value = 2266.3333
$ws.Range("K33").Value = 173.75
$ws.Range("L33").Value = 2266.3333
$ws.Range("M33").Value = 55.25
$ws.Range("N33").Value = -2724.3333

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 10774.6
$ws.Range("J138").Value = 8749.556
$ws.Range("L138").Value = 26248.668
$ws.Range("N138").Value = -36528.66800000001

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3602.625
$ws.Range("I2").Value = 3602.625
$ws.Range("K2").Value = 3602.625
$ws.Range("M2").Value = -3489.625

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3716.25
$ws.Range("I45").Value = 2548.5557
$ws.Range("K45").Value = 2548.5557
$ws.Range("M45").Value = -2171.5557

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1244
$ws.Range("I97").Value = 1143.4445
$ws.Range("K97").Value = 1143.4445
$ws.Range("M97").Value = -647.4445000000001

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 8000
$ws.Range("I110").Value = 8000
$ws.Range("K110").Value = 8000
$ws.Range("M110").Value = -5955

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3602.625
$ws.Range("I116").Value = 3602.625
$ws.Range("K116").Value = 3602.625
$ws.Range("M116").Value = -1308.625

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1675.9375
$ws.Range("I132").Value = 1593.9286
$ws.Range("K132").Value = 4781.7858
$ws.Range("M132").Value = -2251.7858

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3602.625
$ws.Range("I3").Value = 3602.625
$ws.Range("K3").Value = 3602.625
$ws.Range("M3").Value = -3488.625

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3309.3635
$ws.Range("I107").Value = 2721.6667
$ws.Range("K107").Value = 2721.6667
$ws.Range("M107").Value = -801.6667000000002

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1898.5
$ws.Range("I134").Value = 1898.5
$ws.Range("K134").Value = 5695.5
$ws.Range("M134").Value = -3160.5

# CRP row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 17815.385
$ws.Range("I4").Value = 1775
$ws.Range("K4").Value = 1775
$ws.Range("M4").Value = -1663

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 721.36365
$ws.Range("I22").Value = 580
$ws.Range("J22").Value = 1098.3334
$ws.Range("K22").Value = 580
$ws.Range("L22").Value = 1098.3334
$ws.Range("M22").Value = -230
$ws.Range("N22").Value = -1798.3334

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5015.7144
$ws.Range("I31").Value = 4203.3335
$ws.Range("K31").Value = 4203.3335
$ws.Range("M31").Value = -3908.3335

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5015.7144
$ws.Range("I34").Value = 4203.3335
$ws.Range("K34").Value = 4203.3335
$ws.Range("M34").Value = -4001.3335

# CRP row 52
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 50000
$ws.Range("I52").Value = 50000
$ws.Range("K52").Value = 50000
$ws.Range("M52").Value = -49706

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1295.4
$ws.Range("I58").Value = 995
$ws.Range("J58").Value = 1495.6666
$ws.Range("K58").Value = 995
$ws.Range("L58").Value = 1495.6666
$ws.Range("M58").Value = -792
$ws.Range("N58").Value = -1901.6666

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1732
$ws.Range("I134").Value = 1732
$ws.Range("K134").Value = 5196
$ws.Range("M134").Value = -2661

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1295.4
$ws.Range("I136").Value = 995
$ws.Range("J136").Value = 1495.6666
$ws.Range("K136").Value = 2985
$ws.Range("L136").Value = 4486.9998
$ws.Range("M136").Value = -435
$ws.Range("N136").Value = -9586.9998

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 45873.375
$ws.Range("J2").Value = 183433.5
$ws.Range("L2").Value = 1100601
$ws.Range("N2").Value = -1100827

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 702.7143
$ws.Range("I5").Value = 937.25
$ws.Range("J5").Value = 390
$ws.Range("K5").Value = 2811.75
$ws.Range("L5").Value = 1170
$ws.Range("M5").Value = -2699.75
$ws.Range("N5").Value = -1394

# CUL row 21
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 250
$ws.Range("I21").Value = 250
$ws.Range("K21").Value = 750
$ws.Range("M21").Value = -577

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3705.6155
$ws.Range("I68").Value = 4099.143
$ws.Range("J68").Value = 3246.5
$ws.Range("K68").Value = 12297.429
$ws.Range("L68").Value = 9739.5
$ws.Range("M68").Value = -11486.429
$ws.Range("N68").Value = -11361.5

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3705.6155
$ws.Range("I71").Value = 4099.143
$ws.Range("J71").Value = 3246.5
$ws.Range("K71").Value = 36892.287
$ws.Range("L71").Value = 29218.5
$ws.Range("M71").Value = -32836.287
$ws.Range("N71").Value = -37330.5

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1766
$ws.Range("I107").Value = 1803
$ws.Range("J107").Value = 1747.5
$ws.Range("K107").Value = 5409
$ws.Range("L107").Value = 5242.5
$ws.Range("M107").Value = -3489
$ws.Range("N107").Value = -9082.5

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1831.6666
$ws.Range("J132").Value = 1831.6666
$ws.Range("L132").Value = 16484.9994
$ws.Range("N132").Value = -21544.9994

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 615
$ws.Range("I134").Value = 615
$ws.Range("K134").Value = 1845
$ws.Range("M134").Value = 3225

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 702.7143
$ws.Range("I135").Value = 937.25
$ws.Range("J135").Value = 390
$ws.Range("K135").Value = 8435.25
$ws.Range("L135").Value = 3510
$ws.Range("M135").Value = -5900.25
$ws.Range("N135").Value = -8580

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 9500
$ws.Range("I141").Value = 9000
$ws.Range("K141").Value = 27000
$ws.Range("M141").Value = -21820

# GSM row 3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4765.5557
$ws.Range("I3").Value = 4457.647
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 4457.647
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -4341.647
$ws.Range("N3").Value = -10232

# GSM row 4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 3100
$ws.Range("I4").Value = 1200
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = -1088
$ws.Range("N4").Value = -5224

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3271.2856
$ws.Range("I122").Value = 2579.8
$ws.Range("K122").Value = 7739.400000000001
$ws.Range("M122").Value = -5289.400000000001

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1688.8182
$ws.Range("I7").Value = 1928.5
$ws.Range("J7").Value = 1049.6666
$ws.Range("K7").Value = 1928.5
$ws.Range("L7").Value = 1049.6666
$ws.Range("M7").Value = -1816.5
$ws.Range("N7").Value = -1273.6666

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2665.5833
$ws.Range("I22").Value = 2699.2
$ws.Range("K22").Value = 2699.2
$ws.Range("M22").Value = -2404.2

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2665.5833
$ws.Range("I27").Value = 2699.2
$ws.Range("K27").Value = 2699.2
$ws.Range("M27").Value = -2592.2

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1166.6666
$ws.Range("I40").Value = 1166.6666
$ws.Range("K40").Value = 1166.6666
$ws.Range("M40").Value = -1030.6666

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1688.8182
$ws.Range("I126").Value = 1928.5
$ws.Range("J126").Value = 1049.6666
$ws.Range("K126").Value = 5785.5
$ws.Range("L126").Value = 3148.9998
$ws.Range("M126").Value = -3315.5
$ws.Range("N126").Value = -8088.9998

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5003.2085
$ws.Range("I136").Value = 4933.647
$ws.Range("K136").Value = 14800.941
$ws.Range("M136").Value = -12250.941

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2342.5715
$ws.Range("I122").Value = 2483
$ws.Range("K122").Value = 7449
$ws.Range("M122").Value = -4999

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1358.375
$ws.Range("I126").Value = 1123.8572
$ws.Range("K126").Value = 3371.5716
$ws.Range("M126").Value = -901.5715999999998

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2070.5715
$ws.Range("I132").Value = 2082
$ws.Range("J132").Value = 2062
$ws.Range("K132").Value = 6246
$ws.Range("L132").Value = 6186
$ws.Range("M132").Value = -3716
$ws.Range("N132").Value = -11246

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1221.75
$ws.Range("I136").Value = 1221.75
$ws.Range("K136").Value = 3665.25
$ws.Range("M136").Value = -1115.25
